$wb = $excel.ActiveWorkbook

$wsVeda    = $wb.Worksheets.Item("veda input")
$wsThermal = $wb.Worksheets.Item("Thermal elec op costs")
$wsUc      = $wb.Worksheets.Item("uc_data")

# --- "Thermal elec op costs": insert a blank row above row 1 and a blank
#     column to the left of column A (shifts all existing content down/right
#     by one row/column). ---
$wsThermal.Rows.Item(1).Insert()
$wsThermal.Columns.Item(1).Insert()

# --- "uc_data": same shift (insert row above row 1, column before A). ---
$wsUc.Rows.Item(1).Insert()
$wsUc.Columns.Item(1).Insert()

# --- Selection / active-sheet bookkeeping -------------------------------
# "veda input" keeps its own selection but no longer needs a special
# selection set inside it - reset to the top-left cell.
$wsVeda.Activate()
$wsVeda.Range("A1").Select()

# "Thermal elec op costs" loses the tab-selected / special-selection state.
$wsThermal.Activate()
$wsThermal.Range("A1").Select()

# "uc_data" becomes the active / selected tab.
$wsUc.Activate()
$wsUc.Range("A1").Select()

Write-Host "done"
